$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 342) holds a date serial that changes from 45186 to 45188.
$ws.Range("C2:C342").Value = 45188
